$wb = $excel.ActiveWorkbook

# --- Rename "Hoja3" -> "Distancias muros" ---
$ws2 = $wb.Worksheets.Item("Hoja3")
$ws2.Name = "Distancias muros"

# --- Fill in distance values (F6:K32) on "Distancias muros" ---
$ws2.Range("F6").Value = 2.6
$ws2.Range("G6").Value = 0.4
$ws2.Range("H6").Value = 1.4
$ws2.Range("I6").Value = 7.16
$ws2.Range("J6").Value = 5.3
$ws2.Range("K6").Value = 5.34
$ws2.Range("F7").Value = 3.4
$ws2.Range("G7").Value = 0.7
$ws2.Range("H7").Value = 1.5
$ws2.Range("I7").Value = 0.4
$ws2.Range("J7").Value = 2.57
$ws2.Range("K7").Value = 1.52
$ws2.Range("F8").Value = 2
$ws2.Range("G8").Value = 0.7
$ws2.Range("H8").Value = 2.6
$ws2.Range("I8").Value = 0.6
$ws2.Range("J8").Value = 2.57
$ws2.Range("K8").Value = 4.3
$ws2.Range("F9").Value = 0.8
$ws2.Range("G9").Value = 0.7
$ws2.Range("H9").Value = 2
$ws2.Range("I9").Value = 0.6
$ws2.Range("J9").Value = 3.9
$ws2.Range("K9").Value = 1.92
$ws2.Range("F10").Value = 2
$ws2.Range("G10").Value = 0.6
$ws2.Range("H10").Value = 2.8
$ws2.Range("I10").Value = 0.7
$ws2.Range("J10").Value = 1.67
$ws2.Range("K10").Value = 1.6
$ws2.Range("F11").Value = 2.86
$ws2.Range("G11").Value = 0.6
$ws2.Range("H11").Value = 4.4000000000000004
$ws2.Range("I11").Value = 0.7
$ws2.Range("J11").Value = 6.38
$ws2.Range("K11").Value = 2.08
$ws2.Range("F12").Value = 5.35
$ws2.Range("G12").Value = 1.92
$ws2.Range("H12").Value = 0.8
$ws2.Range("I12").Value = 0.7
$ws2.Range("J12").Value = 2.67
$ws2.Range("K12").Value = 5
$ws2.Range("F13").Value = 4
$ws2.Range("G13").Value = 4.8
$ws2.Range("H13").Value = 2
$ws2.Range("I13").Value = 4.9000000000000004
$ws2.Range("J13").Value = 0.85
$ws2.Range("K13").Value = 4.43
$ws2.Range("F14").Value = 2.2999999999999998
$ws2.Range("G14").Value = 3.7
$ws2.Range("H14").Value = 0.8
$ws2.Range("I14").Value = 0.95
$ws2.Range("J14").Value = 0.8
$ws2.Range("K14").Value = 1.92
$ws2.Range("F15").Value = 4.8
$ws2.Range("G15").Value = 1.1599999999999999
$ws2.Range("H15").Value = 2.73
$ws2.Range("I15").Value = 1.92
$ws2.Range("J15").Value = 4
$ws2.Range("K15").Value = 1.92
$ws2.Range("F16").Value = 0.85
$ws2.Range("G16").Value = 5.16
$ws2.Range("H16").Value = 2
$ws2.Range("I16").Value = 4.42
$ws2.Range("J16").Value = 2.2999999999999998
$ws2.Range("K16").Value = 2.6
$ws2.Range("F17").Value = 2
$ws2.Range("G17").Value = 3.84
$ws2.Range("H17").Value = 4
$ws2.Range("I17").Value = 1.2
$ws2.Range("J17").Value = 4.8
$ws2.Range("K17").Value = 3.24
$ws2.Range("F18").Value = 4.4000000000000004
$ws2.Range("G18").Value = 1.18
$ws2.Range("H18").Value = 2.2999999999999998
$ws2.Range("I18").Value = 3.79
$ws2.Range("J18").Value = 2
$ws2.Range("K18").Value = 4.9000000000000004
$ws2.Range("F19").Value = 0.8
$ws2.Range("G19").Value = 6.15
$ws2.Range("H19").Value = 4.8
$ws2.Range("I19").Value = 1.18
$ws2.Range("J19").Value = 8.43
$ws2.Range("K19").Value = 4.92
$ws2.Range("F20").Value = 5.5
$ws2.Range("G20").Value = 4.43
$ws2.Range("H20").Value = 8.4499999999999993
$ws2.Range("I20").Value = 6.1
$ws2.Range("J20").Value = 0.85
$ws2.Range("K20").Value = 6.18
$ws2.Range("F21").Value = 0.62
$ws2.Range("G21").Value = 1.19
$ws2.Range("H21").Value = 0.85
$ws2.Range("I21").Value = 5.14
$ws2.Range("J21").Value = 1.02
$ws2.Range("K21").Value = 3.82
$ws2.Range("F22").Value = 0.62
$ws2.Range("G22").Value = 3.38
$ws2.Range("H22").Value = 1
$ws2.Range("I22").Value = 4.41
$ws2.Range("J22").Value = 2.8
$ws2.Range("K22").Value = 1.2
$ws2.Range("F23").Value = 2.62
$ws2.Range("G23").Value = 0.6
$ws2.Range("H23").Value = 2.75
$ws2.Range("I23").Value = 1.92
$ws2.Range("J23").Value = 5.52
$ws2.Range("K23").Value = 6.15
$ws2.Range("F24").Value = 1.56
$ws2.Range("G24").Value = 0.6
$ws2.Range("H24").Value = 5.5
$ws2.Range("I24").Value = 1.92
$ws2.Range("J24").Value = 0.8
$ws2.Range("K24").Value = 5
$ws2.Range("F25").Value = 1.4
$ws2.Range("G25").Value = 1.92
$ws2.Range("H25").Value = 0.8
$ws2.Range("I25").Value = 0.6
$ws2.Range("J25").Value = 0.8
$ws2.Range("K25").Value = 0.53
$ws2.Range("G26").Value = 1.92
$ws2.Range("H26").Value = 0.5
$ws2.Range("I26").Value = 0.6
$ws2.Range("J26").Value = 0.5
$ws2.Range("K26").Value = 1.77
$ws2.Range("G27").Value = 5.25
$ws2.Range("H27").Value = 0.5
$ws2.Range("I27").Value = 4.83
$ws2.Range("J27").Value = 0.5
$ws2.Range("K27").Value = 8.1
$ws2.Range("G28").Value = 5.03
$ws2.Range("H28").Value = 0.5
$ws2.Range("I28").Value = 3.4
$ws2.Range("J28").Value = 0.5
$ws2.Range("G29").Value = 1.3
$ws2.Range("I29").Value = 5.27
$ws2.Range("G30").Value = 0.4
$ws2.Range("I30").Value = 6.57
$ws2.Range("G31").Value = 2.8
$ws2.Range("G32").Value = 4.92

# --- Comments on "Distancias muros" ---
$c1 = $ws2.Range("I6").AddComment()
$c1.Text("Autor:" + [char]10 + "Muro debajo del estacionamiento piso 2")

$c2 = $ws2.Range("H28").AddComment()
$c2.Text("Autor:" + [char]10 + "Pilares eje 17??")

$c3 = $ws2.Range("J28").AddComment()
$c3.Text("Autor:" + [char]10 + "Pilares eje 17" + [char]10)

# --- Sheet3 selection change (do this before reselecting sheet2, so sheet2 stays the active tab) ---
$ws3 = $wb.Worksheets.Item("Verificación corte muros")
$ws3.Range("J26").Select()

# --- Restore "Distancias muros" as the active sheet/tab and set its selection/scroll ---
$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 5
$ws2.Range("K28:K30").Select()
